$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 862, pushing the existing rows 862..900 down to 864..902
$ws.Rows.Item(862).Insert()
$ws.Rows.Item(862).Insert()

# Populate new row 862
$ws.Range("A862").Value = 7
$ws.Range("B862").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C862").Value = "Ñuble"
$ws.Range("D862").Value = 44939
$ws.Range("E862").Value = 16
$ws.Range("F862").Value = "Fruta"
$ws.Range("G862").Value = 100102
$ws.Range("H862").Value = "Cítricos"
$ws.Range("I862").Value = 100102003
$ws.Range("J862").Value = "Limón"
$ws.Range("K862").Value = "Sin especificar"
$ws.Range("L862").Value = "1a amarillo"
$ws.Range("M862").Value = 160
$ws.Range("N862").Value = 14000
$ws.Range("O862").Value = 15000
$ws.Range("P862").Value = 14500
$ws.Range("Q862").Value = "$/malla 16 kilos"
$ws.Range("R862").Value = "Región de O'Higgins"
$ws.Range("S862").Value = 906
$ws.Range("T862").Value = 16

# Populate new row 863
$ws.Range("A863").Value = 7
$ws.Range("B863").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C863").Value = "Ñuble"
$ws.Range("D863").Value = 44939
$ws.Range("E863").Value = 16
$ws.Range("F863").Value = "Fruta"
$ws.Range("G863").Value = 100102
$ws.Range("H863").Value = "Cítricos"
$ws.Range("I863").Value = 100102003
$ws.Range("J863").Value = "Limón"
$ws.Range("K863").Value = "Sin especificar"
$ws.Range("L863").Value = "1a plateado"
$ws.Range("M863").Value = 160
$ws.Range("N863").Value = 18000
$ws.Range("O863").Value = 19000
$ws.Range("P863").Value = 18500
$ws.Range("Q863").Value = "$/malla 18 kilos"
$ws.Range("R863").Value = "Provincia de Melipilla"
$ws.Range("S863").Value = 1028
$ws.Range("T863").Value = 18

Write-Output "done"
